# FUNCTIONALITY: Wrote two new test cases.
# Update the test-case counters on Sheet1 (B2 = Automated count, C2 = Total count)
# and move the active selection to C2.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 2
$ws.Range("C2").Value = 8

$ws.Range("C2").Select()
